$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("matrix")

# Update age-band labels in column C
$ws.Range("C1").Value = "18-30"
$ws.Range("C2").Value = "31-44"
$ws.Range("C3").Value = "45-56"
$ws.Range("C4").Value = "56-61"
$ws.Range("C5").Value = "62-66"
$ws.Range("C6").Value = "67-75"
$ws.Range("C7").Value = "75above"

# Update premium amounts in column D
$ws.Range("D1").Value = 250
$ws.Range("D2").Value = 500
$ws.Range("D3").Value = 750
$ws.Range("D4").Value = 950
$ws.Range("D5").Value = 1000
$ws.Range("D6").Value = 1100
$ws.Range("D7").Value = 12000

# Move the active selection from A8:D36 to B1
$ws.Range("B1").Select()
